$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-22"

# Update the label for the May row
$ws.Range("A6").Value = "May (through 05-22)"

# Update May row (row 6) values for years 2015-2022 (columns B-I)
$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 44
$ws.Range("E6").Value = 32
$ws.Range("F6").Value = 32
$ws.Range("G6").Value = 42
$ws.Range("H6").Value = 82
$ws.Range("I6").Value = 80

# Update Total row (row 7) values for years 2015-2022 (columns B-I)
$ws.Range("B7").Value = 103
$ws.Range("C7").Value = 195
$ws.Range("D7").Value = 297
$ws.Range("E7").Value = 278
$ws.Range("F7").Value = 187
$ws.Range("G7").Value = 304
$ws.Range("H7").Value = 605
$ws.Range("I7").Value = 632
